# Fixed Elo Decimal Place Inconsistencies
#
# The "Elo Rating" column (and the Home/Away Elo columns in the game
# predictions table) were stored as numbers using a very long decimal
# numeric format ("0.############"). They are corrected to their proper,
# more-precise values and re-entered as text using a simpler "0.##" display
# format. The "Home Team Win Probability" column keeps the long decimal
# format (now on its own style) but is recomputed with slightly different
# precision and is likewise stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- NBA TEAM ELO RATINGS table (column C, rows 5-34) ------------------
$eloRatings = @(
    @{row=5;  val="1797.11"},
    @{row=6;  val="1777.48"},
    @{row=7;  val="1683.89"},
    @{row=8;  val="1661.8"},
    @{row=9;  val="1616.82"},
    @{row=10; val="1612.74"},
    @{row=11; val="1601.3"},
    @{row=12; val="1589.98"},
    @{row=13; val="1588.03"},
    @{row=14; val="1582.82"},
    @{row=15; val="1578.14"},
    @{row=16; val="1551.64"},
    @{row=17; val="1541.99"},
    @{row=18; val="1529.16"},
    @{row=19; val="1502.88"},
    @{row=20; val="1478.74"},
    @{row=21; val="1470.47"},
    @{row=22; val="1465.62"},
    @{row=23; val="1438.25"},
    @{row=24; val="1438.0"},
    @{row=25; val="1436.24"},
    @{row=26; val="1409.99"},
    @{row=27; val="1388.17"},
    @{row=28; val="1378.24"},
    @{row=29; val="1358.1"},
    @{row=30; val="1355.48"},
    @{row=31; val="1325.78"},
    @{row=32; val="1285.63"},
    @{row=33; val="1277.19"},
    @{row=34; val="1262.75"}
)

foreach ($entry in $eloRatings) {
    $cell = $ws.Cells.Item($entry.row, 3)
    $cell.NumberFormat = "@"
    $cell.Value = $entry.val
    $cell.NumberFormat = "0.##"
}

# --- GAME PREDICTIONS table (rows 39-42) --------------------------------
# Home Elo (C), Away Elo (D) re-entered the same way as above, and the
# Home Team Win Probability (F) recomputed with corrected precision.
$games = @(
    @{row=39; c="1797.11"; d="1355.48"; f="0.9576239942153544"},
    @{row=40; c="1578.14"; d="1277.19"; f="0.9095418528145847"},
    @{row=41; c="1551.64"; d="1541.99"; f="0.652760633498995"},
    @{row=42; c="1325.78"; d="1470.47"; f="0.4360382848393158"}
)

foreach ($g in $games) {
    $cCell = $ws.Cells.Item($g.row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $g.c
    $cCell.NumberFormat = "0.##"

    $dCell = $ws.Cells.Item($g.row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $g.d
    $dCell.NumberFormat = "0.##"

    $fCell = $ws.Cells.Item($g.row, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $g.f
    $fCell.NumberFormat = "0.############"
}
